$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1686
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = -1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1266
$ws.Range("L2").Value = 899
$ws.Range("M2").Value = 367
$ws.Range("N2").Value = 242
$ws.Range("O2").Value = 126
$ws.Range("P2").Value = 66
$ws.Range("Q2").Value = 87
$ws.Range("R2").Value = -15
$ws.Range("S2").Value = -68
$ws.Range("T2").Value = 21
$ws.Range("U2").Value = 65
$ws.Range("V2").Value = 432
$ws.Range("W2").Value = 1.17
$ws.Range("X2").Value = 0.01
$ws.Range("Y2").Value = -0.38
$ws.Range("Z2").Value = 0.02
$ws.Range("AA2").Value = 244.8
$ws.Range("AB2").Value = 205.32
$ws.Range("AC2").Value = -14
$ws.Range("AD2").Value = -356.62
$ws.Range("AE2").Value = 3652
$ws.Range("AF2").Value = 1.36
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 1.01
$ws.Range("AI2").Value = -286.33
$ws.Range("AJ2").Value = 6613820

# Row 3
$ws.Range("D3").Value = 1713
$ws.Range("E3").Value = -37
$ws.Range("F3").Value = -33
$ws.Range("G3").Value = -50
$ws.Range("H3").Value = -48
$ws.Range("I3").Value = -44
$ws.Range("J3").Value = -4
$ws.Range("K3").Value = 1282
$ws.Range("L3").Value = 932
$ws.Range("M3").Value = 350
$ws.Range("N3").Value = 210
$ws.Range("O3").Value = 140
$ws.Range("P3").Value = 66
$ws.Range("Q3").Value = 75
$ws.Range("R3").Value = -33
$ws.Range("S3").Value = -51
$ws.Range("T3").Value = 8
$ws.Range("U3").Value = 67
$ws.Range("V3").Value = 403
$ws.Range("W3").Value = -2.15
$ws.Range("X3").Value = -2.77
$ws.Range("Y3").Value = -19.3
$ws.Range("Z3").Value = -3.73
$ws.Range("AA3").Value = 266.08
$ws.Range("AB3").Value = 134.19
$ws.Range("AC3").Value = -659
$ws.Range("AD3").Value = -7.89
$ws.Range("AE3").Value = 3173
$ws.Range("AF3").Value = 1.64
$ws.Range("AG3").Value = 30
$ws.Range("AH3").Value = 0.58
$ws.Range("AI3").Value = -4.55
$ws.Range("AJ3").Value = 6613820

# Row 4
$ws.Range("D4").Value = 1996
$ws.Range("E4").Value = 38
$ws.Range("F4").Value = 38
$ws.Range("G4").Value = 28
$ws.Range("H4").Value = 21
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 14
$ws.Range("K4").Value = 1295
$ws.Range("L4").Value = 938
$ws.Range("M4").Value = 357
$ws.Range("N4").Value = 209
$ws.Range("O4").Value = 148
$ws.Range("P4").Value = 66
$ws.Range("Q4").Value = 75
$ws.Range("R4").Value = 5
$ws.Range("S4").Value = -79
$ws.Range("T4").Value = 14
$ws.Range("U4").Value = 61
$ws.Range("V4").Value = 337
$ws.Range("W4").Value = 1.88
$ws.Range("X4").Value = 1.07
$ws.Range("Y4").Value = 3.52
$ws.Range("Z4").Value = 1.66
$ws.Range("AA4").Value = 262.83
$ws.Range("AB4").Value = 143
$ws.Range("AC4").Value = 111
$ws.Range("AD4").Value = 49.82
$ws.Range("AE4").Value = 3164
$ws.Range("AF4").Value = 1.75
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 0.9
$ws.Range("AI4").Value = 44.89
$ws.Range("AJ4").Value = 6613820

# Row 5
$ws.Range("D5").Value = 2070
$ws.Range("E5").Value = 86
$ws.Range("F5").Value = 86
$ws.Range("G5").Value = 87
$ws.Range("H5").Value = 101
$ws.Range("I5").Value = 86
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = 1341
$ws.Range("L5").Value = 830
$ws.Range("M5").Value = 511
$ws.Range("N5").Value = 358
$ws.Range("O5").Value = 154
$ws.Range("P5").Value = 66
$ws.Range("Q5").Value = 137
$ws.Range("R5").Value = -21
$ws.Range("S5").Value = -72
$ws.Range("T5").Value = 12
$ws.Range("U5").Value = 124
$ws.Range("V5").Value = 265
$ws.Range("W5").Value = 4.14
$ws.Range("X5").Value = 4.88
$ws.Range("Y5").Value = 30.16
$ws.Range("Z5").Value = 7.66
$ws.Range("AA5").Value = 162.26
$ws.Range("AB5").Value = 274.34
$ws.Range("AC5").Value = 1293
$ws.Range("AD5").Value = 14.15
$ws.Range("AE5").Value = 5410
$ws.Range("AF5").Value = 3.38
$ws.Range("AG5").Value = 70
$ws.Range("AH5").Value = 0.38
$ws.Range("AI5").Value = 5.41
$ws.Range("AJ5").Value = 6613820

# Row 6
$ws.Range("D6").Value = 2047
$ws.Range("E6").Value = 70
$ws.Range("F6").Value = 70
$ws.Range("G6").Value = 65
$ws.Range("H6").Value = 70
$ws.Range("I6").Value = 56
$ws.Range("K6").Value = 1282
$ws.Range("L6").Value = 692
$ws.Range("M6").Value = 590
$ws.Range("N6").Value = 424
$ws.Range("P6").Value = 66
$ws.Range("Q6").Value = 99
$ws.Range("R6").Value = -15
$ws.Range("S6").Value = -80
$ws.Range("T6").Value = 26
$ws.Range("U6").Value = 73
$ws.Range("V6").Value = 190
$ws.Range("W6").Value = 3.44
$ws.Range("X6").Value = 3.44
$ws.Range("Y6").Value = 14.29
$ws.Range("Z6").Value = 5.36
$ws.Range("AA6").Value = 117.19
$ws.Range("AB6").Value = 349.42
$ws.Range("AC6").Value = 845
$ws.Range("AD6").Value = 25.68
$ws.Range("AE6").Value = 6416
$ws.Range("AF6").Value = 3.38
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 0.46
$ws.Range("AI6").Value = 11.84
$ws.Range("AJ6").Value = 6613820

# Row 7 - clear data cells, keep A/B/C
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8 - clear data cells, keep A/B/C
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9 - clear data cells, keep A/B/C
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

